$d = $word.ActiveDocument

# 1. Update the generation date/time stamp on the title page.
$d.Content.Find.Execute(
    "May  27, 2021 (11:54:01 PM)", $false, $false, $false, $false, $false,
    $true, 1, $false, "May  28, 2021 (01:53:57 AM)", 2)

# 2. Reword the "latest version" bullet about the textbook availability.
$d.Content.Find.Execute(
    "latest version of textbook is also available under", $false, $false, $false, $false, $false,
    $true, 1, $false, "latest version of built resource is available as a .zip file under", 2)

# 3. Collapse the syntax-highlighted "v1.0.0" example into a single plain run.
$d.Content.Find.Execute(
    "Edit this release, giving it a semantic name and a version, such as v1.0.0. Name and version can be the same.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Edit this release, giving it a semantic name and a version, such as v1.0.0. Name and version can be the same.", 2)

# 4. Wrap "This is a pre-release" in curly quotes and flatten its formatting.
$d.Content.Find.Execute(
    "This is a pre-release", $false, $false, $false, $false, $false,
    $true, 1, $false, [char]0x201C + "This is a pre-release" + [char]0x201D, 2)

# 5. Update the release-distribution sentence.
$d.Content.Find.Execute(
    "university mirror site, box, and archived.", $false, $false, $false, $false, $false,
    $true, 1, $false, "and archived on galileo.", 2)
